# "Update countries & provincias Spain"
#
# The COVID-19 country table is sorted by total cases (column B) descending.
# This refresh brings several rows' stats up to date. A few countries'
# totals now overtake (or tie/re-tie) their neighbours, so those rows swap
# places to keep the table sorted; everything else is untouched.
#
# Each entry below is one worksheet row (by its current row number) with its
# final country name (column A) and, where the stats themselves changed,
# the new Casos totales / Nuevos casos / Casos activos / Recuperados /
# (unused) / Casos criticos / Muertes values for columns B:H.

$updates = @(
    @{ Row = 1;   Name = "Datos actualizados a 26 de Junio de 2020 a las 23:33" },

    @{ Row = 4;   Name = "Estados Unidos";               Vals = @(2545541, 40953, 1059814, 1358386, 0, 561, 127341) },
    @{ Row = 5;   Name = "Brasil";                        Vals = @(1274974, 41827,  673729,  545284, 0, 907,  55961) },
    @{ Row = 10;  Name = "Peru";                          Vals = @( 272364,  3762,  159806,  103619, 0, 178,   8939) },
    @{ Row = 21;  Name = "Sudafrica";                     Vals = @( 124590,  6215,   64111,   58139, 0,  48,   2340) },
    @{ Row = 50;  Name = "Barein";                        Vals = @(  24805,   724,   19137,    5595, 0,   2,     73) },

    # Guayana Francesa overtakes Tailandia
    @{ Row = 96;  Name = "Republica de Africa Central";   Vals = @(   3340,    96,     661,    2639, 0,   0,     40) },
    @{ Row = 97;  Name = "Guayana Francesa";               Vals = @(   3270,   237,    1166,    2093, 0,   1,     11) },
    @{ Row = 98;  Name = "Tailandia";                      Vals = @(   3162,     4,    3040,      64, 0,   0,     58) },

    # Paraguay overtakes Libano and Eslovaquia
    @{ Row = 115; Name = "Paraguay";                       Vals = @(   1711,   142,    1013,     685, 0,   0,     13) },
    @{ Row = 116; Name = "Libano";                         Vals = @(   1697,    35,    1144,     520, 0,   0,     33) },
    @{ Row = 117; Name = "Eslovaquia";                     Vals = @(   1643,    13,    1455,     160, 0,   0,     28) },

    @{ Row = 163; Name = "Siria";                          Vals = @(    255,    13,     102,     145, 0,   1,      8) },
    @{ Row = 176; Name = "Camboya";                        Vals = @(    130,     0,     128,       2, 0,   0,      0) },

    # Tied totals re-sorted (values unchanged, just swap the two countries)
    @{ Row = 200; Name = "Santa Lucia";                    Vals = @(     19,     0,      19,       0, 0,   0,      0) },
    @{ Row = 201; Name = "Laos";                           Vals = @(     19,     0,      19,       0, 0,   0,      0) },
    @{ Row = 202; Name = "Fiyi";                           Vals = @(     18,     0,      18,       0, 0,   0,      0) },
    @{ Row = 203; Name = "Dominica";                       Vals = @(     18,     0,      18,       0, 0,   0,      0) },
    @{ Row = 208; Name = "Groenlandia";                    Vals = @(     13,     0,      13,       0, 0,   0,      0) },
    @{ Row = 209; Name = "Islas Malvinas";                 Vals = @(     13,     0,      13,       0, 0,   0,      0) },
    @{ Row = 212; Name = "Seychelles";                     Vals = @(     11,     0,      11,       0, 0,   0,      0) },
    @{ Row = 213; Name = "Montserrat";                     Vals = @(     11,     0,      10,       0, 0,   0,      1) }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($item in $updates) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Name

    if ($item.ContainsKey("Vals")) {
        $vals = $item.Vals
        for ($i = 0; $i -lt $vals.Length; $i++) {
            $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
        }
    }
}
